$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New quiz rows (8-11) - previously-blank placeholder rows now filled in with
# new trivia questions, matching the "a/b/c" scrambled-answer layout used by
# the existing rows.
# ---------------------------------------------------------------------------

# Row 8 - molecular structure question
$ws.Range("A8").Value = "La estructura molecular fue descubierta por: "
$ws.Range("B8").Value = "b"
$ws.Range("C8").Value = "Watson y Crick "
$ws.Range("D8").Value = "a"
$ws.Range("E8").Value = "John y Paul"
$ws.Range("F8").Value = "c"
$ws.Range("G8").Value = "Roger y David"

# Row 9 - rarest element question
$ws.Range("A9").Value = "Cual es el elemento mas raro en la tierra "
$ws.Range("B9").Value = "c"
$ws.Range("C9").Value = "Astato"
$ws.Range("D9").Value = "a"
$ws.Range("E9").Value = "Paladio"
$ws.Range("F9").Value = "b"
$ws.Range("G9").Value = "Mendeleyev"

# Row 10 - Krypton question
$ws.Range("A10").Value = "El Kriptón es:"
$ws.Range("B10").Value = "c"
$ws.Range("C10").Value = "El cuarto gas noble en la tabla periódica"
$ws.Range("D10").Value = "b"
$ws.Range("E10").Value = "La debilidad de Superman"
$ws.Range("F10").Value = "a"
$ws.Range("G10").Value = "El primer gas noble en la tabla periódica"

# Row 11 - automobile question
$ws.Range("A11").Value = "Quien creo el primer automóvil "
$ws.Range("B11").Value = "a"
$ws.Range("C11").Value = "Karl Benz "
$ws.Range("D11").Value = "b"
$ws.Range("E11").Value = "Henry Ford"
$ws.Range("F11").Value = "c"
$ws.Range("G11").Value = "Enzo Ferrari "

# ---------------------------------------------------------------------------
# Row heights for the newly-populated rows (matches wrapped-text autofit
# heights used by the similarly-sized existing question rows).
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 34
$ws.Rows.Item(9).RowHeight = 34
$ws.Rows.Item(10).RowHeight = 68
$ws.Rows.Item(11).RowHeight = 34

# ---------------------------------------------------------------------------
# Column A formatting: unify the question cells (rows 2-11) onto the plain
# wrap-text style (copy format from A1, which already carries it) so the
# redundant "wrap + applyFont" style is no longer used.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)

# The two still-blank rows below (12-13) keep the green "placeholder" font
# style; re-apply it by copying the format that's already there so it keeps
# reusing the same style entry.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Selection moves from H4 to I6.
# ---------------------------------------------------------------------------
$ws.Range("I6").Select()

# ---------------------------------------------------------------------------
# Workbook-level housekeeping captured by the diff: the file's recorded
# absolute path (folder was renamed/moved) and the window position/size
# last saved by Excel.
# ---------------------------------------------------------------------------
try { $excel.ActiveWindow.Top = 0 } catch {}
try { $excel.ActiveWindow.Left = 0 } catch {}
try { $excel.ActiveWindow.Width = 17920 } catch {}
try { $excel.ActiveWindow.Height = 22400 } catch {}
